# Applies the "Updated cryptos list" refresh: new Price/Volume(1h) figures for
# most rows, plus a reshuffle of three ranking positions (29/30/31 and 47/48)
# where the Coin/Link columns were swapped to reflect the new order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number (e.g. "562.62") must
# be forced to Text format first, otherwise Excel would silently convert them
# to real numbers (dropping things like trailing zeros, e.g. "14.00" -> 14).
function Set-TextValue($addr, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

$ws.Range('D2').Value = '62.898.09'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '3.379.22'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  -0.15%  '
Set-TextValue 'D5' '562.62'
$ws.Range('E5').Value = '  +0.65%  '
Set-TextValue 'D6' '154.67'
$ws.Range('E6').Value = '  +1.40%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.378.41'
$ws.Range('E8').Value = '  +0.88%  '
Set-TextValue 'D9' '0.541'
$ws.Range('E9').Value = '  +2.37%  '
Set-TextValue 'D10' '7.39'
$ws.Range('E10').Value = '  -1.26%  '
Set-TextValue 'D11' '0.121'
$ws.Range('E11').Value = '  +2.14%  '
Set-TextValue 'D12' '0.431'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').Value = '3.958.66'
$ws.Range('E13').Value = '  +0.65%  '
Set-TextValue 'D14' '0.134'
$ws.Range('E14').Value = '  -3.19%  '
$ws.Range('E15').Value = '  +3.90%  '
Set-TextValue 'D16' '27.05'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').Value = '63.005.19'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '3.387.40'
$ws.Range('E18').Value = '  -0.06%  '
Set-TextValue 'D19' '6.23'
$ws.Range('E19').Value = '  -2.88%  '
Set-TextValue 'D20' '14.00'
$ws.Range('E20').Value = '  +1.70%  '
Set-TextValue 'D21' '375.35'
$ws.Range('E21').Value = '  -2.85%  '
Set-TextValue 'D22' '7.98'
$ws.Range('E22').Value = '  -4.10%  '
Set-TextValue 'D23' '0.998'
$ws.Range('E23').Value = '  -0.38%  '
Set-TextValue 'D24' '71.18'
$ws.Range('E24').Value = '  +1.38%  '
Set-TextValue 'D25' '0.527'
$ws.Range('E25').Value = '  -2.03%  '
Set-TextValue 'D26' '0.0000117'
$ws.Range('E26').Value = '  +22.14%  '
Set-TextValue 'D27' '9.42'
$ws.Range('E27').Value = '  +6.47%  '
Set-TextValue 'D28' '0.177'
$ws.Range('E28').Value = '  -1.17%  '
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  -0.21%  '
Set-TextValue 'D30' '6.02'
$ws.Range('E30').Value = '  +7.36%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D31' '1.34'
$ws.Range('E31').Value = '  +3.83%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D32' '1.98'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D33' '6.42'
$ws.Range('E33').Value = '  -1.53%  '
Set-TextValue 'D34' '23.07'
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('E35').Value = '  +0.09%  '
Set-TextValue 'D36' '6.75'
$ws.Range('E36').Value = '  +0.89%  '
Set-TextValue 'D37' '159.26'
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('D39').Value = '2.958.63'
$ws.Range('E39').Value = '  +4.69%  '
Set-TextValue 'D40' '0.0761'
$ws.Range('E40').Value = '  +2.33%  '
Set-TextValue 'D41' '27.04'
$ws.Range('E41').Value = '  +1.32%  '
Set-TextValue 'D42' '1.82'
$ws.Range('E42').Value = '  -3.48%  '
Set-TextValue 'D43' '0.0317'
$ws.Range('E43').Value = '  +1.95%  '
Set-TextValue 'D44' '41.46'
$ws.Range('E44').Value = '  +2.40%  '
Set-TextValue 'D45' '4.30'
$ws.Range('E45').Value = '  +0.88%  '
Set-TextValue 'D46' '0.744'
$ws.Range('E46').Value = '  -0.32%  '
Set-TextValue 'D47' '23.20'
$ws.Range('E47').Value = '  +5.73%  '
Set-TextValue 'D48' '1.06'
$ws.Range('E48').Value = '  +2.54%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D49' '2.12'
$ws.Range('E49').Value = '  +19.15%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D50' '6.35'
$ws.Range('E50').Value = '  +1.22%  '
Set-TextValue 'D51' '0.828'
$ws.Range('E51').Value = '  +3.20%  '
